$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.738.46"
$ws.Range("D3").Value = "1.891.27"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.79"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4758"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2940"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06546"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.12"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07739"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7424"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.94"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").Value = "1.887.91"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.256"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "276.28"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "30.805.94"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.23"
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007577"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "2.134.11"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.342"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.255"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.260"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.20"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.87"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.928"
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09742"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.313"
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.184"
$ws.Range("E33").Value = "  +2.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04897"
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.129"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7012"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01915"
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.352"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.87"
$ws.Range("E41").Value = "  +6.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.035"
$ws.Range("E42").Value = "  +3.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4267"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8443"
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9996"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.61"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.401"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.088"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.72"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "920.21"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05775"
$ws.Range("E51").Value = "  +2.27%  "
